$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 253.45
$ws.Range("I38").Value = 24.6
$ws.Range("J38").Value = 940
$ws.Range("K38").Value = 73.80000000000001
$ws.Range("L38").Value = 2820
$ws.Range("M38").Value = 298.2
$ws.Range("N38").Value = -3564

# Row 125
$ws.Range("H125").Value = 8305.736999999999
$ws.Range("I125").Value = 14948.777
$ws.Range("J125").Value = 2327
$ws.Range("K125").Value = 134538.993
$ws.Range("L125").Value = 20943
$ws.Range("M125").Value = -132078.993
$ws.Range("N125").Value = -25863

# Row 137
$ws.Range("H137").Value = 4359.189
$ws.Range("I137").Value = 2841.2341
$ws.Range("J137").Value = 6018.3486
$ws.Range("K137").Value = 8523.702300000001
$ws.Range("L137").Value = 18055.0458
$ws.Range("M137").Value = -5973.702300000001
$ws.Range("N137").Value = -23155.0458

# Row 138
$ws.Range("H138").Value = 4085.6785
$ws.Range("I138").Value = 1679.2
$ws.Range("J138").Value = 6862.385
$ws.Range("K138").Value = 5037.6
$ws.Range("L138").Value = 20587.155
$ws.Range("M138").Value = 102.3999999999996
$ws.Range("N138").Value = -30867.155

# Row 139
$ws.Range("H139").Value = 46995
$ws.Range("J139").Value = 46995
$ws.Range("L139").Value = 46995
$ws.Range("N139").Value = -57275

# Row 140
$ws.Range("H140").Value = 43983.332
$ws.Range("J140").Value = 43983.332
$ws.Range("L140").Value = 43983.332
$ws.Range("N140").Value = -54343.332

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1990
$ws.Range("I107").Value = 1105.5
$ws.Range("J107").Value = 2579.6667
$ws.Range("K107").Value = 1105.5
$ws.Range("L107").Value = 2579.6667
$ws.Range("M107").Value = 814.5
$ws.Range("N107").Value = -6419.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 3879.6765
$ws.Range("I62").Value = 4211.2144
$ws.Range("K62").Value = 4211.2144
$ws.Range("M62").Value = -3587.2144

# Row 65
$ws.Range("H65").Value = 3879.6765
$ws.Range("I65").Value = 4211.2144
$ws.Range("K65").Value = 21056.072
$ws.Range("M65").Value = -17936.072

# Row 132
$ws.Range("H132").Value = 38468244
$ws.Range("I132").Value = 62508700
$ws.Range("J132").Value = 3517.4
$ws.Range("K132").Value = 187526100
$ws.Range("L132").Value = 10552.2
$ws.Range("M132").Value = -187523570
$ws.Range("N132").Value = -15612.2

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 756.1111
$ws.Range("I5").Value = 514.8461
$ws.Range("J5").Value = 980.1429000000001
$ws.Range("K5").Value = 1544.5383
$ws.Range("L5").Value = 2940.4287
$ws.Range("M5").Value = -1432.5383
$ws.Range("N5").Value = -3164.4287

# Row 33
$ws.Range("H33").Value = 78.97143
$ws.Range("I33").Value = 78.51613
$ws.Range("J33").Value = 82.5
$ws.Range("K33").Value = 471.09678
$ws.Range("L33").Value = 495
$ws.Range("M33").Value = -188.09678
$ws.Range("N33").Value = -1061

# Row 38
$ws.Range("H38").Value = 681.6
$ws.Range("I38").Value = 31.09091
$ws.Range("J38").Value = 1192.7142
$ws.Range("K38").Value = 93.27273
$ws.Range("L38").Value = 3578.1426
$ws.Range("M38").Value = 253.72727
$ws.Range("N38").Value = -4272.142599999999

# Row 39
$ws.Range("H39").Value = 3060
$ws.Range("J39").Value = 3750
$ws.Range("L39").Value = 11250
$ws.Range("N39").Value = -11838

# Row 40
$ws.Range("H40").Value = 360.9091
$ws.Range("I40").Value = 95
$ws.Range("J40").Value = 680
$ws.Range("K40").Value = 380
$ws.Range("L40").Value = 2720
$ws.Range("M40").Value = -311
$ws.Range("N40").Value = -2858

# Row 41
$ws.Range("H41").Value = 361.66666
$ws.Range("I41").Value = 300
$ws.Range("J41").Value = 423.33334
$ws.Range("K41").Value = 900
$ws.Range("L41").Value = 1270.00002
$ws.Range("M41").Value = -562
$ws.Range("N41").Value = -1946.00002

# Row 42
$ws.Range("H42").Value = 2225
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 2225
$ws.Range("K42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("M42").Value = 6675
$ws.Range("N42").Value = -7743

# Row 43
$ws.Range("H43").Value = 3000
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 9000
$ws.Range("N43").Value = -9228

# Row 44
$ws.Range("H44").Value = 333333340
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

# Row 46
$ws.Range("H46").Value = 950
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -1109
$ws.Range("N46").Value = -4682

# Row 47
$ws.Range("H47").Value = 700
$ws.Range("I47").Value = 100
$ws.Range("J47").Value = 850
$ws.Range("K47").Value = 300
$ws.Range("L47").Value = 2550
$ws.Range("M47").Value = 131
$ws.Range("N47").Value = -3412

# Row 48
$ws.Range("H48").Value = 2840
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 2840
$ws.Range("K48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("M48").Value = 8520
$ws.Range("N48").Value = -9020

# Row 49
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1000
$ws.Range("K49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("M49").Value = 3000
$ws.Range("N49").Value = -3312

# Row 51
$ws.Range("H51").Value = 525
$ws.Range("I51").Value = 525
$ws.Range("K51").Value = 1575
$ws.Range("M51").Value = -1115

# Row 54
$ws.Range("H54").Value = 5000
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16118

# Row 55
$ws.Range("H55").Value = 2695.4783
$ws.Range("I55").Value = 1440
$ws.Range("J55").Value = 3044.2222
$ws.Range("K55").Value = 4320
$ws.Range("L55").Value = 9132.6666
$ws.Range("M55").Value = -4143
$ws.Range("N55").Value = -9486.6666

# Row 68
$ws.Range("H68").Value = 861.05884
$ws.Range("I68").Value = 636.8148
$ws.Range("J68").Value = 1113.3334
$ws.Range("K68").Value = 1910.4444
$ws.Range("L68").Value = 3340.0002
$ws.Range("M68").Value = -1099.4444
$ws.Range("N68").Value = -4962.0002

# Row 71
$ws.Range("H71").Value = 861.05884
$ws.Range("I71").Value = 636.8148
$ws.Range("J71").Value = 1113.3334
$ws.Range("K71").Value = 5731.3332
$ws.Range("L71").Value = 10020.0006
$ws.Range("M71").Value = -1675.3332
$ws.Range("N71").Value = -18132.0006

# Row 107
$ws.Range("H107").Value = 521.4286
$ws.Range("I107").Value = 238
$ws.Range("J107").Value = 1726
$ws.Range("K107").Value = 714
$ws.Range("L107").Value = 5178
$ws.Range("M107").Value = 1206
$ws.Range("N107").Value = -9018

# Row 135
$ws.Range("H135").Value = 756.1111
$ws.Range("I135").Value = 514.8461
$ws.Range("J135").Value = 980.1429000000001
$ws.Range("K135").Value = 4633.6149
$ws.Range("L135").Value = 8821.286100000001
$ws.Range("M135").Value = -2098.6149
$ws.Range("N135").Value = -13891.2861

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 3367.7273
$ws.Range("I100").Value = 1268.4
$ws.Range("J100").Value = 5117.1665
$ws.Range("K100").Value = 1268.4
$ws.Range("L100").Value = 5117.1665
$ws.Range("M100").Value = -727.4000000000001
$ws.Range("N100").Value = -6199.1665

# Row 136
$ws.Range("H136").Value = 5004.1377
$ws.Range("I136").Value = 1321.6316
$ws.Range("J136").Value = 12000.9
$ws.Range("K136").Value = 3964.8948
$ws.Range("L136").Value = 36002.7
$ws.Range("M136").Value = -1414.8948
$ws.Range("N136").Value = -41102.7

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 2447.3684
$ws.Range("I62").Value = 2555.0908
$ws.Range("J62").Value = 2299.25
$ws.Range("K62").Value = 2555.0908
$ws.Range("L62").Value = 2299.25
$ws.Range("M62").Value = -1931.0908
$ws.Range("N62").Value = -3547.25

# Row 65
$ws.Range("H65").Value = 2447.3684
$ws.Range("I65").Value = 2555.0908
$ws.Range("J65").Value = 2299.25
$ws.Range("K65").Value = 12775.454
$ws.Range("L65").Value = 11496.25
$ws.Range("M65").Value = -9655.454
$ws.Range("N65").Value = -17736.25
